$wb = $excel.ActiveWorkbook

# Sheet at position 3 (rId3) currently holds the SanDisk SSD row under name "ID_c4955e1".
# Sheet at position 4 (rId4) currently holds the John Hardy bracelet row under name "ID_2b6fb5b".
$s3 = $wb.Worksheets.Item(3)
$s4 = $wb.Worksheets.Item(4)

# Rename sheet 4 away first so its name is free for sheet 3 to take (avoids a name collision),
# and give it the brand-new id used for the newly tracked product.
$s4.Name = "ID_9037e8e"
$s4.Cells.Item(2, 1).Value = 109.95
$s4.Cells.Item(2, 3).Value = "Fjallraven - Foldsack No. 1 Backpack, Fits 15 Laptops"

# Sheet 3 now takes over the identity/content that used to live on sheet 4 (the John Hardy product).
$s3.Name = "ID_2b6fb5b"
$s3.Cells.Item(2, 1).Value = 695
$s3.Cells.Item(2, 3).Value = "John Hardy Women's Legends Naga Gold & Silver Dragon Station Chain Bracelet"
